$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '70.223.53'
$ws.Range('E2').Value = '  -1.12%  '
$ws.Range('D3').Value = '3.506.08'
$ws.Range('E3').Value = '  -1.28%  '
$ws.Range('E4').Value = '  +0.11%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '604.56'
$ws.Range('E5').Value = '  -0.52%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '172.82'
$ws.Range('E6').Value = '  -1.38%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.608'
$ws.Range('E7').Value = '  -1.64%  '
$ws.Range('D8').Value = '3.499.99'
$ws.Range('E8').Value = '  -1.36%  '
$ws.Range('E9').Value = '  -0.01%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '0.195'
$ws.Range('E10').Value = '  -3.02%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '7.22'
$ws.Range('E11').Value = '  +6.89%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '0.586'
$ws.Range('E12').Value = '  -0.38%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '46.18'
$ws.Range('E13').Value = '  -3.49%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '0.0000276'
$ws.Range('E14').Value = '  -2.02%  '
$ws.Range('D15').Value = '4.079.51'
$ws.Range('E15').Value = '  -1.11%  '
$ws.Range('E16').Value = '  -1.23%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '612.76'
$ws.Range('E17').Value = '  -2.90%  '
$ws.Range('D18').Value = '3.507.00'
$ws.Range('E18').Value = '  -1.25%  '
$ws.Range('D19').Value = '70.234.57'
$ws.Range('E19').Value = '  -0.97%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '17.50'
$ws.Range('E21').Value = '  -0.04%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '0.878'
$ws.Range('E22').Value = '  -1.52%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '9.13'
$ws.Range('E23').Value = '  -9.79%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '98.63'
$ws.Range('E24').Value = '  +1.35%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '15.62'
$ws.Range('E25').Value = '  -2.32%  '
$ws.Range('E26').Value = '  -3.99%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '0.999'
$ws.Range('E27').Value = '  -0.09%  '
$ws.Range('E28').Value = '  -2.69%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '33.92'
$ws.Range('E29').Value = '  +1.16%  '
$ws.Range('E31').Value = '  -5.02%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '8.05'
$ws.Range('E32').Value = '  -5.31%  '
$ws.Range('E33').Value = '  -5.23%  '
$ws.Range('B34').Value = 'Bittensor'
$ws.Range('C34').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '631.64'
$ws.Range('E34').Value = '  +10.71%  '
$ws.Range('B35').Value = 'NEARProtocol'
$ws.Range('C35').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '6.83'
$ws.Range('E35').Value = '  -3.81%  '
$ws.Range('E36').Value = '  -2.84%  '
$ws.Range('E37').Value = '  -0.75%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '0.0483'
$ws.Range('E38').Value = '  +6.28%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '3.48'
$ws.Range('E39').Value = '  -5.72%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '56.83'
$ws.Range('E40').Value = '  -1.29%  '
$ws.Range('E41').Value = '  +0.20%  '
$ws.Range('E42').Value = '  +0.91%  '
$ws.Range('D43').Value = '3.363.02'
$ws.Range('E43').Value = '  +0.29%  '
$ws.Range('D44').Value = '0.0₃0734'
$ws.Range('E44').Value = '  +1.42%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '0.310'
$ws.Range('E45').Value = '  -6.20%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '2.91'
$ws.Range('E46').Value = '  -4.65%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '31.88'
$ws.Range('E47').Value = '  -4.30%  '
$ws.Range('E48').Value = '  -4.69%  '
$ws.Range('E49').Value = '  +0.10%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '133.16'
$ws.Range('E50').Value = '  -0.91%  '
$ws.Range('E51').Value = '  -0.02%  '

$ws.Range('D5').Style = "Normal"
$ws.Range('D6').Style = "Normal"
$ws.Range('D7').Style = "Normal"
$ws.Range('D10').Style = "Normal"
$ws.Range('D11').Style = "Normal"
$ws.Range('D12').Style = "Normal"
$ws.Range('D13').Style = "Normal"
$ws.Range('D14').Style = "Normal"
$ws.Range('D17').Style = "Normal"
$ws.Range('D21').Style = "Normal"
$ws.Range('D22').Style = "Normal"
$ws.Range('D23').Style = "Normal"
$ws.Range('D24').Style = "Normal"
$ws.Range('D25').Style = "Normal"
$ws.Range('D27').Style = "Normal"
$ws.Range('D29').Style = "Normal"
$ws.Range('D32').Style = "Normal"
$ws.Range('D34').Style = "Normal"
$ws.Range('D35').Style = "Normal"
$ws.Range('D38').Style = "Normal"
$ws.Range('D39').Style = "Normal"
$ws.Range('D40').Style = "Normal"
$ws.Range('D45').Style = "Normal"
$ws.Range('D46').Style = "Normal"
$ws.Range('D47').Style = "Normal"
$ws.Range('D50').Style = "Normal"
